$wb = $excel.ActiveWorkbook

# 1. Rename sheets 3,4,5 to prefix with "$表名"
$ws3 = $wb.Worksheets.Item("=rare#别名例子")
$ws3.Name = "`$表名=rare#别名例子"

$ws4 = $wb.Worksheets.Item("+subList#list子表")
$ws4.Name = "`$表名+subList#list子表"

$ws5 = $wb.Worksheets.Item("+subMap#map子表")
$ws5.Name = "`$表名+subMap#map子表"

# 2. Edit sheet1 ("title前六行") cells
$ws1 = $wb.Worksheets.Item("title前六行")
$ws1.Range("G7").Value = "紫"
$ws1.Range("N7").Value = "a|b|c"
$ws1.Range("O7").Value = "a|1&b|2&c|3"
$ws1.Range("G8").Value = "表名(可省略)=别名例子"

# 3. Page setup (print setup dialog touched on sheet1)
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# 4. Selections left by the author while editing
$ws1.Range("K25").Select() | Out-Null
$ws5.Range("B41").Select() | Out-Null
